$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "How many curves can I load in one go?"
$ws.Range("B4").Value = "deepseek1.5"
$ws.Range("C4").Value = "You can load up to 450 curves at a time."
